$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "bands" table (D4:E7) got reversed (band1..band4 -> band4..band1) ---
$ws.Range("D4").Value = "band4"
$ws.Range("E4").Value = "256-192"
$ws.Range("D5").Value = "band3"
$ws.Range("E5").Value = "256-128"
$ws.Range("D6").Value = "band2"
$ws.Range("E6").Value = "256-64"
$ws.Range("D7").Value = "band1"
$ws.Range("E7").Value = "256-0"

# --- Start integrating "preview" fields: give preview_file_start / preview_file_dur real values ---
$ws.Range("B12").Value = 50
$ws.Range("B13").Value = 10

# --- New "stringmap_to" data column (H) gets its own width + light formatting ---
$ws.Columns.Item(8).ColumnWidth = 11.52
$ws.Range("H1").Borders.Item(9).LineStyle = 1
$ws.Range("H3").Borders.Item(9).LineStyle = 1
$ws.Range("H4").Borders.Item(9).LineStyle = 1

# --- Selection moves down to the next empty row of the preview block ---
$ws.Range("B14").Select() | Out-Null
